$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.033.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.54%  "
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.851.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.971.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.408.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.178"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "326.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0517"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.577"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.404"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0224"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("E51").Value = "  -1.19%  "
